# Trade #9 closed at 2026-02-17 12:27:28 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.04
$summary.Range("B4").Value = 0.04
$summary.Range("B5").Value = 0.09
$summary.Range("B6").Value = 9
$summary.Range("B7").Value = 4
$summary.Range("B9").Value = 44.44

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.04
$status.Range("D4").Value = 9
$status.Range("E4").Value = 0.04
$status.Range("F4").Value = 0.04
$status.Range("G4").Value = 44.44

# --- New trade row data, appended as row 10 on both "All Trades" and "MarketMaking" sheets ---
$newRow = @(9, "2026-02-17", "12:27:22", "MarketMaking", "DOWN", 0.88, 0.91, "CLOSED", 3.4091, 0.03, 100.04, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $val = $newRow[$i]
        if ($col -eq 2) {
            # B column holds a "yyyy-mm-dd" looking date that must stay literal
            # text (matching the rest of the column) rather than being
            # auto-converted to a date serial number by COM's smart entry.
            $ws.Cells.Item(10, $col).Value = "'" + $val
            $ws.Cells.Item(10, $col).Style = "Normal"
        } else {
            $ws.Cells.Item(10, $col).Value = $val
        }
    }
}
